# Translate the student-template header row (and the bound table's column
# headers, which stay in sync automatically) from Arabic to English, turn
# off the table's AutoFilter, and move the active selection from D3 to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "gender"
$ws.Range("C1").Value = "birthdate"
$ws.Range("D1").Value = "university_number"
$ws.Range("E1").Value = "study_type"
$ws.Range("F1").Value = "study_year"
$ws.Range("G1").Value = "program"
$ws.Range("H1").Value = "phone"
$ws.Range("I1").Value = "email"
$ws.Range("J1").Value = "notes"

$lo = $ws.ListObjects.Item(1)
$lo.ShowAutoFilter = $false

$ws.Range("A3").Select()
